$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 366.77777
$ws.Range("I18").Value = 245.25
$ws.Range("J18").Value = 464
$ws.Range("K18").Value = 245.25
$ws.Range("L18").Value = 464
$ws.Range("M18").Value = 38.75
$ws.Range("N18").Value = -1032

$ws.Range("H112").Value = 25001700
$ws.Range("J112").Value = 1755.7894
$ws.Range("L112").Value = 5267.3682
$ws.Range("N112").Value = -7483.3682

$ws.Range("H129").Value = 832.0517
$ws.Range("J129").Value = 969.6445
$ws.Range("L129").Value = 2908.9335
$ws.Range("N129").Value = -12908.9335

$ws.Range("H132").Value = 47627780
$ws.Range("I132").Value = 58832376
$ws.Range("K132").Value = 176497128
$ws.Range("M132").Value = -176494598

$ws.Range("H137").Value = 2553.4583
$ws.Range("I137").Value = 1216.36
$ws.Range("J137").Value = 4006.8262
$ws.Range("K137").Value = 3649.08
$ws.Range("L137").Value = 12020.4786
$ws.Range("M137").Value = -1099.08
$ws.Range("N137").Value = -17120.4786

$ws.Range("H138").Value = 5013.862
$ws.Range("I138").Value = 1091.1765
$ws.Range("J138").Value = 5879.909
$ws.Range("K138").Value = 3273.5295
$ws.Range("L138").Value = 17639.727
$ws.Range("M138").Value = 1866.4705
$ws.Range("N138").Value = -27919.727


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3863.1177
$ws.Range("I32").Value = 3529.6206
$ws.Range("K32").Value = 3529.6206
$ws.Range("M32").Value = -3242.6206

$ws.Range("H61").Value = 1093.5636
$ws.Range("I61").Value = 827.8913
$ws.Range("K61").Value = 827.8913
$ws.Range("M61").Value = -615.8913

$ws.Range("H74").Value = 3409.8948
$ws.Range("I74").Value = 3366.7097
$ws.Range("K74").Value = 3366.7097
$ws.Range("M74").Value = -2492.7097

$ws.Range("H77").Value = 3409.8948
$ws.Range("I77").Value = 3366.7097
$ws.Range("K77").Value = 16833.5485
$ws.Range("M77").Value = -12465.5485

$ws.Range("H122").Value = 2047.5807
$ws.Range("I122").Value = 1403.9524
$ws.Range("K122").Value = 4211.857199999999
$ws.Range("M122").Value = -1761.857199999999

$ws.Range("H132").Value = 2488.3022
$ws.Range("I132").Value = 1486.5518
$ws.Range("J132").Value = 4563.357
$ws.Range("K132").Value = 4459.6554
$ws.Range("L132").Value = 13690.071
$ws.Range("M132").Value = -1929.6554
$ws.Range("N132").Value = -18750.071

$ws.Range("H136").Value = 1093.5636
$ws.Range("I136").Value = 827.8913
$ws.Range("K136").Value = 2483.6739
$ws.Range("M136").Value = 66.32610000000022


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1662.9214
$ws.Range("I134").Value = 1013.13434
$ws.Range("J134").Value = 3641.818
$ws.Range("K134").Value = 3039.40302
$ws.Range("L134").Value = 10925.454
$ws.Range("M134").Value = -504.4030199999997
$ws.Range("N134").Value = -15995.454


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 542.9355
$ws.Range("I22").Value = 337.47827
$ws.Range("J22").Value = 1133.625
$ws.Range("K22").Value = 337.47827
$ws.Range("L22").Value = 1133.625
$ws.Range("M22").Value = 12.52172999999999
$ws.Range("N22").Value = -1833.625

$ws.Range("H31").Value = 6495707.5
$ws.Range("I31").Value = 1202.8909
$ws.Range("J31").Value = 22731970
$ws.Range("K31").Value = 1202.8909
$ws.Range("L31").Value = 22731970
$ws.Range("M31").Value = -907.8909000000001
$ws.Range("N31").Value = -22732560

$ws.Range("H34").Value = 6495707.5
$ws.Range("I34").Value = 1202.8909
$ws.Range("J34").Value = 22731970
$ws.Range("K34").Value = 1202.8909
$ws.Range("L34").Value = 22731970
$ws.Range("M34").Value = -1000.8909
$ws.Range("N34").Value = -22732374

$ws.Range("H58").Value = 1767.6046
$ws.Range("I58").Value = 1709.5593
$ws.Range("J58").Value = 1894.4445
$ws.Range("K58").Value = 1709.5593
$ws.Range("L58").Value = 1894.4445
$ws.Range("M58").Value = -1506.5593
$ws.Range("N58").Value = -2300.4445

$ws.Range("H99").Value = 20007288
$ws.Range("I99").Value = 40004376
$ws.Range("J99").Value = 10200
$ws.Range("K99").Value = 40004376
$ws.Range("L99").Value = 10200
$ws.Range("M99").Value = -40002878
$ws.Range("N99").Value = -13196

$ws.Range("H105").Value = 2209.75
$ws.Range("I105").Value = 2031.6666
$ws.Range("J105").Value = 2744
$ws.Range("K105").Value = 2031.6666
$ws.Range("L105").Value = 2744
$ws.Range("M105").Value = -284.6666
$ws.Range("N105").Value = -6238

$ws.Range("H126").Value = 20007288
$ws.Range("I126").Value = 40004376
$ws.Range("J126").Value = 10200
$ws.Range("K126").Value = 120013128
$ws.Range("L126").Value = 30600
$ws.Range("M126").Value = -120010658
$ws.Range("N126").Value = -35540

$ws.Range("H132").Value = 2573.4102
$ws.Range("I132").Value = 2371.64
$ws.Range("J132").Value = 2933.7144
$ws.Range("K132").Value = 7114.92
$ws.Range("L132").Value = 8801.143199999999
$ws.Range("M132").Value = -4584.92
$ws.Range("N132").Value = -13861.1432

$ws.Range("H134").Value = 3655.7551
$ws.Range("I134").Value = 4259.724
$ws.Range("K134").Value = 12779.172
$ws.Range("M134").Value = -10244.172

$ws.Range("H136").Value = 1767.6046
$ws.Range("I136").Value = 1709.5593
$ws.Range("J136").Value = 1894.4445
$ws.Range("K136").Value = 5128.6779
$ws.Range("L136").Value = 5683.333500000001
$ws.Range("M136").Value = -2578.6779
$ws.Range("N136").Value = -10783.3335


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 787.92554
$ws.Range("I131").Value = 467.5
$ws.Range("J131").Value = 826.0714
$ws.Range("K131").Value = 1402.5
$ws.Range("L131").Value = 2478.2142
$ws.Range("M131").Value = 3637.5
$ws.Range("N131").Value = -12558.2142


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2665.5
$ws.Range("I122").Value = 1461.6666
$ws.Range("J122").Value = 5245.143
$ws.Range("K122").Value = 4384.9998
$ws.Range("L122").Value = 15735.429
$ws.Range("M122").Value = -1934.9998
$ws.Range("N122").Value = -20635.429

$ws.Range("H132").Value = 2416.5417
$ws.Range("I132").Value = 1523.909
$ws.Range("J132").Value = 4380.3335
$ws.Range("K132").Value = 4571.727000000001
$ws.Range("L132").Value = 13141.0005
$ws.Range("M132").Value = -2041.727000000001
$ws.Range("N132").Value = -18201.0005


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6260.6665
$ws.Range("I122").Value = 2976.25
$ws.Range("K122").Value = 8928.75
$ws.Range("M122").Value = -6478.75

$ws.Range("H132").Value = 3353.726
$ws.Range("I132").Value = 1194.7709
$ws.Range("J132").Value = 7498.92
$ws.Range("K132").Value = 3584.3127
$ws.Range("L132").Value = 22496.76
$ws.Range("M132").Value = -1054.3127
$ws.Range("N132").Value = -27556.76

$ws.Range("H136").Value = 1979.125
$ws.Range("I136").Value = 1090.0212
$ws.Range("K136").Value = 3270.063599999999
$ws.Range("M136").Value = -720.0635999999995


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3174.1365
$ws.Range("I122").Value = 1839.4375
$ws.Range("J122").Value = 6733.3335
$ws.Range("K122").Value = 5518.3125
$ws.Range("L122").Value = 20200.0005
$ws.Range("M122").Value = -3068.3125
$ws.Range("N122").Value = -25100.0005

$ws.Range("H123").Value = 36194
$ws.Range("J123").Value = 36194
$ws.Range("L123").Value = 36194
$ws.Range("N123").Value = -45994

$ws.Range("H132").Value = 6062062.5
$ws.Range("I132").Value = 616.6111
$ws.Range("J132").Value = 17546908
$ws.Range("K132").Value = 1849.8333
$ws.Range("L132").Value = 52640724
$ws.Range("M132").Value = 680.1667000000002
$ws.Range("N132").Value = -52645784

$ws.Range("H136").Value = 2034.2458
$ws.Range("I136").Value = 747.2308
$ws.Range("J136").Value = 4315.773
$ws.Range("K136").Value = 2241.6924
$ws.Range("L136").Value = 12947.319
$ws.Range("M136").Value = 308.3076000000001
$ws.Range("N136").Value = -18047.319

